$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Exam date: 23 Desember 2025 -> 16 Desember 2025 (both occurrences)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("23 Desember 2025", $false, $false, $false, $false, $false,
                         $true, 1, $false, "16 Desember 2025", 2)

# ---------------------------------------------------------------------------
# 2. Exam start time: 10:00 -> 13:00
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("10:00 ", $false, $false, $false, $false, $false,
                         $true, 1, $false, "13:00 ", 2)

# ---------------------------------------------------------------------------
# 3. Examiner name: Sindy Devila, S.Si, M.Si. -> Prof. Alhadi Bustamam.,
#    S.Si., M.Kom., Ph.D. (both occurrences).
#
#    One of the two occurrences is immediately preceded, in the same
#    paragraph, by a separate run containing just "." whose run formatting
#    (rPr) is identical to the name run's formatting. A plain text
#    replacement on the name run would make Word's run-merging logic fold
#    that "." run into the replaced run (losing the run boundary that the
#    target markup preserves). To avoid that, we briefly nudge the
#    formatting of the immediately preceding character so it no longer
#    matches, perform the text swap, then restore the original formatting.
# ---------------------------------------------------------------------------
$oldName = "Sindy Devila, S.Si, M.Si."
$newName = "Prof. Alhadi Bustamam., S.Si., M.Kom., Ph.D."

$searchRange = $d.Content
$searchRange.Find.Text = $oldName
$searchRange.Find.Forward = $true
$searchRange.Find.Wrap = 0
$searchRange.Find.MatchCase = $false
$searchRange.Find.MatchWholeWord = $false

$found = $searchRange.Find.Execute()
while ($found) {
    $matchStart = $searchRange.Start
    $matchEnd = $searchRange.End

    $precededByRun = $matchStart -gt 0
    if ($precededByRun) {
        $guard = $d.Range($matchStart - 1, $matchStart)
        $guardWasBold = $guard.Bold
        if ($guardWasBold -eq 1) {
            $guard.Bold = 0
        } else {
            $guard.Bold = 1
        }
    }

    $target = $d.Range($matchStart, $matchEnd)
    $target.Text = $newName

    if ($precededByRun) {
        $guard2 = $d.Range($matchStart - 1, $matchStart)
        $guard2.Bold = $guardWasBold
    }

    $searchRange = $d.Range($matchStart + $newName.Length, $d.Content.End)
    $searchRange.Find.Text = $oldName
    $searchRange.Find.Forward = $true
    $searchRange.Find.Wrap = 0
    $searchRange.Find.MatchCase = $false
    $searchRange.Find.MatchWholeWord = $false
    $found = $searchRange.Find.Execute()
}
